$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.028.62"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.238.25"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.73%  "
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  -5.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.08%  "
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "2.579.91"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "2.236.73"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.826"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "43.925.96"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "0.0₃0961"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.61%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("E38").Value = "  -11.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.60%  "
$ws.Range("E41").Value = "  -8.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.46%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "1.737.67"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "85.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.17%  "
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.35%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.49%  "
